$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix three misspelled Car IDs (letter "O" -> digit "0")
$ws.Range("A7").Value = "FD06FCS006"
$ws.Range("A15").Value = "GM09CMR014"
$ws.Range("A38").Value = "HO05ODY037"

# Add new computed columns F (2-digit year code) and G (car age)
$ws.Range("F2").Formula = "=MID(A2,3,2)"
$ws.Range("G2").Formula = "=IF(25-F2<0,100-F2+25,25-F2)"

$ws.Range("F3:F53").Formula = "=MID(A3,3,2)"
$ws.Range("G3:G53").Formula = "=IF(25-F3<0,100-F3+25,25-F3)"

# Update the view: scroll back to top and move selection to C38
[void]$ws.Range("C38").Select()
